# Changed "Attribut" to "Attribute" for consistency, and fixed the row
# ordering in the "ARC" sheet for IfcSpace (Name/LongName/PredefinedType/
# IsExternal/IsInteriorOrExteriorSpace/LongName-test block) and the
# IfcWindow/IfcDoor PredefinedType rows.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Header row: rename "Attribut..." columns to "Attribute..." for consistency.
    if ($ws.Cells.Item(1,7).Value2 -eq "AttributDescriptionIT") {
        $ws.Range("G1").Value = "AttributeDescriptionIT"
    }
    if ($ws.Cells.Item(1,8).Value2 -eq "AttributName") {
        $ws.Range("H1").Value = "AttributeName"
    }
}

# Fix the row order on the "ARC" sheet (rows 8-11 and 12-15).
$arc = $wb.Worksheets.Item("ARC")

$arc.Range("H8").Value = "PredefinedType"
$arc.Range("K8").Value = "EXTERNAL, INTERNAL"

$arc.Range("F9").Value = "Pset_SpaceCommon"
$arc.Range("H9").Value = "IsInteriorOrExteriorSpace"
$arc.Range("K9").Value = "TRUE, FALSE"

$arc.Range("H10").Value = "IsExternal"
$arc.Range("K10").Value = "TRUE, FALSE"

$arc.Range("F11").Value = ""
$arc.Range("H11").Value = "LongName"
$arc.Range("K11").Value = "test"

$arc.Range("E12").Value = "IfcWindow"
$arc.Range("K12").Value = ""

$arc.Range("E13").Value = "IfcWindow"

$arc.Range("E14").Value = "IfcDoor"
$arc.Range("K14").Value = "DOOR, GATE, TRAPDOOR, USERDEFINED, NOTDEFINED"

$arc.Range("E15").Value = "IfcDoor"
